$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two username/mid values (stored as text) in column A
$ws.Range("A2").Value = "20150914015"
$ws.Range("A3").Value = "20150914016"

# Move the active selection from A3 to B4
$ws.Range("B4").Select()
